$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.768.01'
$ws.Range('E2').Value = '  +1.13%  '
$ws.Range('D3').Value = '2.273.36'
$ws.Range('E3').Value = '  +0.15%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '121.31'
$ws.Range('E5').Value = '  +6.95%  '
$ws.Range('E6').Value = '  +0.85%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.650'
$ws.Range('E7').Value = '  +4.87%  '
$ws.Range('E8').Value = '  +0.26%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.628'
$ws.Range('E9').Value = '  +5.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '48.44'
$ws.Range('E10').Value = '  +0.47%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0950'
$ws.Range('E11').Value = '  +2.67%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '9.26'
$ws.Range('E12').Value = '  +6.03%  '
$ws.Range('E13').Value = '  -0.52%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.66'
$ws.Range('E14').Value = '  +1.76%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.915'
$ws.Range('E15').Value = '  +6.92%  '
$ws.Range('D16').Value = '2.617.19'
$ws.Range('E16').Value = '  +0.24%  '
$ws.Range('D17').Value = '2.274.58'
$ws.Range('E17').Value = '  +0.15%  '
$ws.Range('D18').Value = '43.700.87'
$ws.Range('E18').Value = '  +1.25%  '
$ws.Range('E19').Value = '  +3.18%  '
$ws.Range('E20').Value = '  +0.92%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.33'
$ws.Range('E21').Value = '  +1.52%  '
$ws.Range('E22').Value = '  -0.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.86'
$ws.Range('E23').Value = '  +2.26%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.60'
$ws.Range('E24').Value = '  -2.08%  '
$ws.Range('E25').Value = '  +1.91%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.97'
$ws.Range('E26').Value = '  +5.73%  '
$ws.Range('E27').Value = '  +1.79%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '43.25'
$ws.Range('E28').Value = '  +6.18%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.42'
$ws.Range('E29').Value = '  +1.59%  '
$ws.Range('E30').Value = '  +0.49%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '173.62'
$ws.Range('E31').Value = '  +1.45%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.70'
$ws.Range('E32').Value = '  +2.09%  '
$ws.Range('E33').Value = '  +2.39%  '
$ws.Range('E34').Value = '  +2.81%  '
$ws.Range('E35').Value = '  +3.93%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.28'
$ws.Range('E36').Value = '  +11.60%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0385'
$ws.Range('E37').Value = '  +9.83%  '
$ws.Range('E38').Value = '  +0.40%  '
$ws.Range('E39').Value = '  +4.45%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.55'
$ws.Range('E40').Value = '  +5.10%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '74.09'
$ws.Range('E41').Value = '  -2.18%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '13.73'
$ws.Range('E42').Value = '  -3.81%  '
$ws.Range('E43').Value = '  +2.39%  '
$ws.Range('E44').Value = '  +0.06%  '
$ws.Range('B45').Value = 'ARBITRUM'
$ws.Range('C45').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.40'
$ws.Range('E45').Value = '  +2.30%  '
$ws.Range('B46').Value = 'THORChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.92'
$ws.Range('E46').Value = '  -3.01%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '74.94'
$ws.Range('E47').Value = '  +42.86%  '
$ws.Range('E48').Value = '  +4.10%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '103.27'
$ws.Range('E49').Value = '  +2.70%  '
$ws.Range('E50').Value = '  +2.18%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.48'
$ws.Range('E51').Value = '  -1.56%  '
